$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Class_EU for the two "small / mechanical" sludge-separator rows
# from "primary" to "secondary"
$ws.Range("E14").Value = "secondary"
$ws.Range("E15").Value = "secondary"

# Move the active selection to E10 (matches the final saved selection)
$ws.Range("E10").Select()
